# Auto-generated edit script: updates crypto price/volume figures
# to match the refreshed GitHub Actions data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.369.79'
$ws.Range('E2').Value = '  +2.67%  '
$ws.Range('D3').Value = '1.845.56'
$ws.Range('E3').Value = '  +2.34%  '
$ws.Range('E4').Value = '  +0.30%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '230.12'
$r.ClearFormats()
$ws.Range('E5').Value = '  +2.48%  '
$ws.Range('E6').Value = '  +3.83%  '
$ws.Range('E7').Value = '  +0.17%  '
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '42.69'
$r.ClearFormats()
$ws.Range('E8').Value = '  +12.24%  '
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '0.307'
$r.ClearFormats()
$ws.Range('E9').Value = '  +7.09%  '
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.0694'
$r.ClearFormats()
$ws.Range('E10').Value = '  +3.61%  '
$ws.Range('E11').Value = '  +3.59%  '
$ws.Range('D12').Value = '2.113.71'
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('D13').Value = '1.844.35'
$ws.Range('E13').Value = '  +2.22%  '
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '11.35'
$r.ClearFormats()
$ws.Range('E14').Value = '  +2.99%  '
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '0.672'
$r.ClearFormats()
$ws.Range('E15').Value = '  +7.44%  '
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '4.67'
$r.ClearFormats()
$ws.Range('E16').Value = '  +6.84%  '
$ws.Range('D17').Value = '35.352.57'
$ws.Range('E17').Value = '  +2.74%  '
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '70.36'
$r.ClearFormats()
$ws.Range('E18').Value = '  +3.74%  '
$ws.Range('D19').Value = '0.0₃0801'
$ws.Range('E19').Value = '  +4.32%  '
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '245.94'
$r.ClearFormats()
$ws.Range('E20').Value = '  +1.94%  '
$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '12.11'
$r.ClearFormats()
$ws.Range('E21').Value = '  +10.07%  '
$ws.Range('E22').Value = '  +13.87%  '
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('E24').Value = '  +0.43%  '
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '169.21'
$r.ClearFormats()
$ws.Range('E25').Value = '  -0.65%  '
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '7.91'
$r.ClearFormats()
$ws.Range('E26').Value = '  +3.11%  '
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '17.77'
$r.ClearFormats()
$ws.Range('E27').Value = '  +1.92%  '
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '1.39'
$r.ClearFormats()
$ws.Range('E29').Value = '  +13.63%  '
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').Value = '3.361.44'
$ws.Range('E31').Value = '  +38.35%  '
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '0.0546'
$r.ClearFormats()
$ws.Range('E32').Value = '  +6.62%  '
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '3.94'
$r.ClearFormats()
$ws.Range('E33').Value = '  +5.13%  '
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '4.06'
$r.ClearFormats()
$ws.Range('E34').Value = '  +5.63%  '
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '1.86'
$r.ClearFormats()
$ws.Range('E35').Value = '  +2.83%  '
$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '97.08'
$r.ClearFormats()
$ws.Range('E36').Value = '  +19.20%  '
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '0.689'
$r.ClearFormats()
$ws.Range('E37').Value = '  +7.94%  '
$ws.Range('D38').Value = '1.348.89'
$ws.Range('E38').Value = '  +2.08%  '
$ws.Range('E39').Value = '  +3.29%  '
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '2.44'
$r.ClearFormats()
$ws.Range('E40').Value = '  +7.01%  '
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '0.0195'
$r.ClearFormats()
$ws.Range('E41').Value = '  +3.70%  '
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '1.00'
$r.ClearFormats()
$ws.Range('E42').Value = '  +6.41%  '
$ws.Range('E43').Value = '  +3.83%  '
$ws.Range('E44').Value = '  +7.53%  '
$ws.Range('E45').Value = '  +0.64%  '
$ws.Range('E46').Value = '  +0.30%  '
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '6.23'
$r.ClearFormats()
$ws.Range('E47').Value = '  +8.79%  '
$ws.Range('E48').Value = '  +2.11%  '
$ws.Range('D49').Value = '2.013.46'
$ws.Range('E49').Value = '  +2.48%  '
$ws.Range('E50').Value = '  +0.22%  '
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '103.58'
$r.ClearFormats()
$ws.Range('E51').Value = '  +1.73%  '
